$d = $word.ActiveDocument

# Change 1: "CS_08_02_CO" -> "CS_09_02_CO" (guion reference number)
$d.Content.Find.Execute("CS_08_02_CO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CS_09_02_CO", 2)

# Change 2: "Refuerza tu aprendizaje: Los caudillos" -> "Refuerza tu aprendizaje: los caudillos"
$d.Content.Find.Execute("Refuerza tu aprendizaje: Los caudillos", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Refuerza tu aprendizaje: los caudillos", 2)

# Change 3: "independencia,América" -> "independencia, América" (missing space).
# Replace the whole "caudillo, independencia,América Latina" span (but not the
# lone-space run right before it) so the run -- and the stray proofErr
# spell/grammar-check markers inside it -- collapse into a single plain run,
# matching the reverted document.
$d.Content.Find.Execute("caudillo, independencia,América Latina", $true, $false, $false, $false, $false,
                         $true, 1, $false, "caudillo, independencia, América Latina", 2)

# Change 4: "...pueda validarlas" -> "...pueda validarlas." (trailing period)
$d.Content.Find.Execute("pueda validarlas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "pueda validarlas.", 2)

# Change 5: drop the word "tres " from "Describe por lo menos tres características"
$d.Content.Find.Execute("Describe por lo menos tres características", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Describe por lo menos características", 2)

# Change 5b: add trailing space after the final sentence of that same paragraph
$d.Content.Find.Execute("Independencia. Explica cada una de ellas.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Independencia. Explica cada una de ellas. ", 2)

# The text replacements above removed the "_GoBack" bookmark that previously sat
# between "Refuerza tu aprendizaje: L" and "os caudillos". Word relocates this
# bookmark to mark the last edited spot -- put it back right after "recordar",
# matching the reverted document's layout.
$findRng = $d.Content
$findRng.Find.Execute("recordar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($findRng.End, $findRng.End)
$bmRange.Bookmarks.Add("_GoBack")
